# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (same per-fund layout/styling as the
# existing "2021-Q3"/"2021-Q4" sheets) right before the "总计" (totals)
# sheet, and prepends a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" worksheet by duplicating "2021-Q4" (same
#    headers/column layout/styling) and dropping it in right before
#    "总计" - this keeps header/row formatting identical to the other
#    quarterly sheets.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# "2021-Q4" has 3 fund rows (010343, 000049, 010344); 2022-Q1 only has
# 2 (010343, 010344), so remove the middle "000049" row - this leaves
# the 010344 row (originally row 4) as row 3, with its fund code/name
# already correct.
$q1.Rows.Item(3).Delete()

# Fix the leading index column (0-based row counter) after the delete.
$q1.Range("A3").Value = 1

# Columns D:G hold figures stored as TEXT (not numbers) in the source
# data - force text formatting first so values like "0.20"/"93.65"
# keep their exact original digits instead of being parsed as numbers.
$q1.Range("D2:G3").NumberFormat = "@"

# Row 2 - 华宝英国富时100指数（QDII）A : update the figures that changed
# for 2022-Q1 (code/name/index are already correct from the template).
$q1.Range("D2").Value = "0.20"
$q1.Range("E2").Value = "93.65"
$q1.Range("F2").Value = "3.24"
$q1.Range("G2").Value = "0.0065"
$q1.Range("H2").Value = 9

# Row 3 - 华宝英国富时100指数（QDII）C : update the figures that changed
# for 2022-Q1.
$q1.Range("D3").Value = "0.06"
$q1.Range("E3").Value = "93.65"
$q1.Range("F3").Value = "3.24"
$q1.Range("G3").Value = "0.0019"
$q1.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" (totals) sheet, pushing the
#    existing rows down by one and renumbering the leading index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing rows (previously r=2 "2021-Q4", r=3 "2021-Q3") shift down to
# r=3 / r=4. Rewrite them bottom-up so we never clobber unread data.
# A4 is a brand-new cell (row 4 didn't exist before), so copy A3's
# style (bordered/centered index-column look) onto it first.
$total.Range("A3").Copy($total.Range("A4"))
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.01

# New row for 2022-Q1 at the top of the data (r=2)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01
